# Generate Report for Handback
# Adds a "version mismatch" handback-validation result for the
# 4bac05a3-1c7a-4008-aa32-c45c9c8d6837 entry on both the zh-cn and de-de
# report sheets: a hyperlinked "Latest Target File" value, the generated
# handback xlf name, the handback datetime, and the error detail message.

$wb = $excel.ActiveWorkbook

# Column width bump (OOXML stores width = ColumnWidth + 5/6); use that
# offset so the saved <col width="..."> comes out to exactly 40.
$targetOoxmlWidth = 40
$columnWidthValue = $targetOoxmlWidth - (5 / 6)

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/910fdc9bfab63b6e5a6c3570b3d71a5465b2d6c8/e2e/4bac05a3-1c7a-4008-aa32-c45c9c8d6837.md"
$displayMd = "4bac05a3-1c7a-4008-aa32-c45c9c8d6837.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24755287abe191ec19dd3c85165068ca04f735e1/e2e/4bac05a3-1c7a-4008-aa32-c45c9c8d6837.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/910fdc9bfab63b6e5a6c3570b3d71a5465b2d6c8/e2e/4bac05a3-1c7a-4008-aa32-c45c9c8d6837.md."

$sheetsInfo = @(
    @{ Name = "zh-cn"; HandbackFile = "4bac05a3-1c7a-4008-aa32-c45c9c8d6837.eceaddae60e5dfacb494baff3098ee3b37c4d5a4.zh-cn.xlf"; HandbackDate = "2016-10-20 00:01:14" },
    @{ Name = "de-de"; HandbackFile = "4bac05a3-1c7a-4008-aa32-c45c9c8d6837.eceaddae60e5dfacb494baff3098ee3b37c4d5a4.de-de.xlf"; HandbackDate = "2016-10-20 00:01:33" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the "Latest Target File", "Latest Handback File" and
    # "Error Detail" columns (I, J, P) to fit the new long values.
    $ws.Columns.Item(9).ColumnWidth = $columnWidthValue
    $ws.Columns.Item(10).ColumnWidth = $columnWidthValue
    $ws.Columns.Item(16).ColumnWidth = $columnWidthValue

    # Row 5 corresponds to the 4bac05a3-1c7a-4008-aa32-c45c9c8d6837 file.
    $ws.Range("J5").Value = $info.HandbackFile
    $ws.Range("K5").Value = $info.HandbackDate
    $ws.Range("P5").Value = $errorDetail

    # I5 ("Latest Target File") becomes a hyperlink to the source .md file,
    # matching the style used by the other hyperlinked cells in column A.
    $ws.Hyperlinks.Add($ws.Range("I5"), $targetUrl, "", "", $displayMd)
    $ws.Range("I5").Font.Name = "Calibri"
    $ws.Range("I5").Font.Underline = 2
    $ws.Range("I5").Font.Color = 15570276
}
